# chore: update Sheets via scheduled runner
#
# Refreshes cached FFXIV Market Board prices (currentAveragePrice* /
# LevePrice*/LeveProfit* columns) across all eight job sheets, and strips
# the bold+bordered header-row formatting (Table headers revert to the
# default "Normal" style) that a prior run had introduced.

$wb = $excel.ActiveWorkbook

# ---- 1) Header style reset -------------------------------------------------
# Row 1 (the Table header) on every sheet loses its bold font + thin border;
# cells revert to the workbook default ("Normal") style.
foreach ($ws in $wb.Worksheets) {
    $ws.Range("A1:N1").Style = "Normal"
}

# ---- 2) Per-cell market-price value updates --------------------------------
# Columns: H currentAveragePrice, I currentAveragePriceNQ, J currentAveragePriceHQ,
#          K LevePriceNQ, L LevePriceHQ, M LeveProfitNQ, N LeveProfitHQ
# A handful of rows also gain/lose an M or N cell entirely when a price
# category becomes (un)available, matching the source diff exactly.

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 90997
$ws.Range("J3").Value = 90997
$ws.Range("L3").Value = 90997
$ws.Range("N3").Value = -91225
$ws.Range("H43").Value = 4671.2144
$ws.Range("J43").Value = 5439.8
$ws.Range("L43").Value = 5439.8
$ws.Range("N43").Value = -5577.8
$ws.Range("H74").Value = 6248.4165
$ws.Range("I74").Value = 4999.8887
$ws.Range("K74").Value = 4999.8887
$ws.Range("M74").Value = -4063.8887
$ws.Range("H77").Value = 6248.4165
$ws.Range("I77").Value = 4999.8887
$ws.Range("K77").Value = 24999.4435
$ws.Range("M77").Value = -20319.4435
$ws.Range("H98").Value = 1576.6
$ws.Range("I98").Value = 970
$ws.Range("J98").Value = 4003
$ws.Range("K98").Value = 970
$ws.Range("L98").Value = 4003
$ws.Range("M98").Value = 528
$ws.Range("N98").Value = -6999
$ws.Range("H102").Value = 90997
$ws.Range("J102").Value = 90997
$ws.Range("L102").Value = 90997
$ws.Range("N102").Value = -97487
$ws.Range("H106").Value = 17100.953
$ws.Range("I106").Value = 14239
$ws.Range("K106").Value = 14239
$ws.Range("M106").Value = -13608
$ws.Range("H122").Value = 1576.6
$ws.Range("I122").Value = 970
$ws.Range("J122").Value = 4003
$ws.Range("K122").Value = 2910
$ws.Range("L122").Value = 12009
$ws.Range("M122").Value = -460
$ws.Range("N122").Value = -16909
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null   # cell removed
$ws.Range("H135").Value = 2055.5
$ws.Range("I135").Value = 1630.3636
$ws.Range("K135").Value = 14673.2724
$ws.Range("M135").Value = -12138.2724
$ws.Range("H137").Value = 4540.0293
$ws.Range("I137").Value = 4463.4
$ws.Range("K137").Value = 13390.2
$ws.Range("M137").Value = -10840.2
$ws.Range("H138").Value = 3793.4
$ws.Range("J138").Value = 6504.3335
$ws.Range("L138").Value = 19513.0005
$ws.Range("N138").Value = -29793.0005

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 4014
$ws.Range("I36").Value = 3281.6
$ws.Range("K36").Value = 3281.6
$ws.Range("M36").Value = -2935.6
$ws.Range("H61").Value = 5928.909
$ws.Range("I61").Value = 2522
$ws.Range("K61").Value = 2522
$ws.Range("M61").Value = -2310
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = $null   # cell removed
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 1899.75
$ws.Range("I132").Value = 1899.75
$ws.Range("K132").Value = 5699.25
$ws.Range("M132").Value = -3169.25
$ws.Range("H136").Value = 5928.909
$ws.Range("I136").Value = 2522
$ws.Range("K136").Value = 7566
$ws.Range("M136").Value = -5016

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1091.5
$ws.Range("I94").Value = 1181
$ws.Range("J94").Value = 554.5
$ws.Range("K94").Value = 1181
$ws.Range("L94").Value = 554.5
$ws.Range("M94").Value = -730
$ws.Range("N94").Value = -1456.5
$ws.Range("H106").Value = 50000
$ws.Range("J106").Value = 50000
$ws.Range("L106").Value = 50000
$ws.Range("N106").Value = -52524   # cell added
$ws.Range("H107").Value = 3982.1
$ws.Range("I107").Value = 3813.4443
$ws.Range("K107").Value = 3813.4443
$ws.Range("M107").Value = -1893.4443
$ws.Range("H134").Value = 2767
$ws.Range("I134").Value = 1874
$ws.Range("J134").Value = 4999.5
$ws.Range("K134").Value = 5622
$ws.Range("L134").Value = 14998.5
$ws.Range("M134").Value = -3087
$ws.Range("N134").Value = -20068.5

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 385.6
$ws.Range("I7").Value = 356.125
$ws.Range("J7").Value = 700
$ws.Range("K7").Value = 356.125
$ws.Range("L7").Value = 700
$ws.Range("M7").Value = -243.125
$ws.Range("N7").Value = -926
$ws.Range("H22").Value = 49712.96
$ws.Range("I22").Value = 54879.363
$ws.Range("J22").Value = 11826
$ws.Range("K22").Value = 54879.363
$ws.Range("L22").Value = 11826
$ws.Range("M22").Value = -54529.363
$ws.Range("N22").Value = -12526
$ws.Range("H28").Value = 12162
$ws.Range("J28").Value = 12162
$ws.Range("L28").Value = 12162
$ws.Range("N28").Value = -12652
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = $null   # cell removed
$ws.Range("H58").Value = 3304.0386
$ws.Range("J58").Value = 4375.8335
$ws.Range("L58").Value = 4375.8335
$ws.Range("N58").Value = -4781.8335
$ws.Range("H70").Value = 67499.5
$ws.Range("J70").Value = 67499.5
$ws.Range("L70").Value = 67499.5
$ws.Range("N70").Value = -68129.5
$ws.Range("H73").Value = 67499.5
$ws.Range("J73").Value = 67499.5
$ws.Range("L73").Value = 67499.5
$ws.Range("N73").Value = -69683.5
$ws.Range("H86").Value = 33700.875
$ws.Range("I86").Value = 8337
$ws.Range("J86").Value = 48919.2
$ws.Range("K86").Value = 8337
$ws.Range("L86").Value = 48919.2
$ws.Range("M86").Value = -7214
$ws.Range("N86").Value = -51165.2
$ws.Range("H89").Value = 33700.875
$ws.Range("I89").Value = 8337
$ws.Range("J89").Value = 48919.2
$ws.Range("K89").Value = 41685
$ws.Range("L89").Value = 244596
$ws.Range("M89").Value = -36069
$ws.Range("N89").Value = -255828
$ws.Range("H95").Value = 30252.6
$ws.Range("J95").Value = 30252.6
$ws.Range("L95").Value = 30252.6
$ws.Range("N95").Value = -35744.6
$ws.Range("H99").Value = 7632.5
$ws.Range("I99").Value = 9332.667
$ws.Range("J99").Value = 5932.3335
$ws.Range("K99").Value = 9332.667
$ws.Range("L99").Value = 5932.3335
$ws.Range("M99").Value = -7834.666999999999
$ws.Range("N99").Value = -8928.3335
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = $null   # cell removed
$ws.Range("H107").Value = 693.5625
$ws.Range("I107").Value = 590.36365
$ws.Range("K107").Value = 590.36365
$ws.Range("M107").Value = 1329.63635
$ws.Range("H126").Value = 7632.5
$ws.Range("I126").Value = 9332.667
$ws.Range("J126").Value = 5932.3335
$ws.Range("K126").Value = 27998.001
$ws.Range("L126").Value = 17797.0005
$ws.Range("M126").Value = -25528.001
$ws.Range("N126").Value = -22737.0005
$ws.Range("H136").Value = 3304.0386
$ws.Range("J136").Value = 4375.8335
$ws.Range("L136").Value = 13127.5005
$ws.Range("N136").Value = -18227.5005

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 128.57143
$ws.Range("J23").Value = 155
$ws.Range("L23").Value = 465
$ws.Range("N23").Value = -935
$ws.Range("H34").Value = 8066.1665
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 8066.1665
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 24198.4995
$ws.Range("M34").Value = $null   # cell removed
$ws.Range("N34").Value = -24366.4995
$ws.Range("H55").Value = 8183.1113
$ws.Range("J55").Value = 7956.125
$ws.Range("L55").Value = 23868.375
$ws.Range("N55").Value = -24222.375

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").Value = $null   # cell removed
$ws.Range("H46").Value = 12571.143
$ws.Range("I46").Value = 9000
$ws.Range("J46").Value = 13999.6
$ws.Range("K46").Value = 9000
$ws.Range("L46").Value = 13999.6
$ws.Range("M46").Value = -8844
$ws.Range("N46").Value = -14311.6
$ws.Range("H68").Value = 58366.11
$ws.Range("J68").Value = 59411.875
$ws.Range("L68").Value = 59411.875
$ws.Range("N68").Value = -61033.875
$ws.Range("H71").Value = 58366.11
$ws.Range("J71").Value = 59411.875
$ws.Range("L71").Value = 178235.625
$ws.Range("N71").Value = -186347.625
$ws.Range("H75").Value = 49499.5
$ws.Range("J75").Value = 49499.5
$ws.Range("L75").Value = 49499.5
$ws.Range("N75").Value = -51247.5
$ws.Range("H78").Value = 49499.5
$ws.Range("J78").Value = 49499.5
$ws.Range("L78").Value = 148498.5
$ws.Range("N78").Value = -157234.5
$ws.Range("H105").Value = 17335.5
$ws.Range("J105").Value = 17335.5
$ws.Range("L105").Value = 17335.5
$ws.Range("N105").Value = -24323.5
$ws.Range("H132").Value = 2829.6
$ws.Range("I132").Value = 2877.3333
$ws.Range("K132").Value = 8631.999899999999
$ws.Range("M132").Value = -6101.999899999999

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 583.3333
$ws.Range("I22").Value = 525
$ws.Range("K22").Value = 525
$ws.Range("M22").Value = -230
$ws.Range("H27").Value = 583.3333
$ws.Range("I27").Value = 525
$ws.Range("K27").Value = 525
$ws.Range("M27").Value = -418
$ws.Range("H40").Value = 15777.667
$ws.Range("I40").Value = 20166.5
$ws.Range("J40").Value = 7000
$ws.Range("K40").Value = 20166.5
$ws.Range("L40").Value = 7000
$ws.Range("M40").Value = -20030.5
$ws.Range("N40").Value = -7272
$ws.Range("H46").Value = 5600.684
$ws.Range("I46").Value = 2025
$ws.Range("J46").Value = 6554.2
$ws.Range("K46").Value = 2025
$ws.Range("L46").Value = 6554.2
$ws.Range("M46").Value = -1837
$ws.Range("N46").Value = -6930.2
$ws.Range("H122").Value = 22249
$ws.Range("I122").Value = 22249
$ws.Range("K122").Value = 66747
$ws.Range("M122").Value = -64297
$ws.Range("H132").Value = 22308.75
$ws.Range("I132").Value = 13366.75
$ws.Range("K132").Value = 40100.25
$ws.Range("M132").Value = -37570.25
$ws.Range("H136").Value = 17649
$ws.Range("I136").Value = 12599.25
$ws.Range("J136").Value = 27748.5
$ws.Range("K136").Value = 37797.75
$ws.Range("L136").Value = 83245.5
$ws.Range("M136").Value = -35247.75
$ws.Range("N136").Value = -88345.5

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 336333.34
$ws.Range("I32").Value = 336333.34
$ws.Range("K32").Value = 336333.34
$ws.Range("M32").Value = -336016.34
$ws.Range("H54").Value = 19999.666
$ws.Range("J54").Value = 19999.666
$ws.Range("L54").Value = 19999.666
$ws.Range("N54").Value = -21039.666
$ws.Range("H58").Value = 20998.727
$ws.Range("I58").Value = 13998.333
$ws.Range("J58").Value = 29399.2
$ws.Range("K58").Value = 13998.333
$ws.Range("L58").Value = 29399.2
$ws.Range("M58").Value = -13690.333
$ws.Range("N58").Value = -30015.2
$ws.Range("H81").Value = 5543.5557
$ws.Range("I81").Value = 2841.7144
$ws.Range("J81").Value = 15000
$ws.Range("K81").Value = 5683.4288
$ws.Range("L81").Value = 30000
$ws.Range("M81").Value = -4622.4288
$ws.Range("N81").Value = -32122
$ws.Range("H84").Value = 5543.5557
$ws.Range("I84").Value = 2841.7144
$ws.Range("J84").Value = 15000
$ws.Range("K84").Value = 28417.144
$ws.Range("L84").Value = 150000
$ws.Range("M84").Value = -23113.144
$ws.Range("N84").Value = -160608
$ws.Range("H100").Value = 1124.625
$ws.Range("I100").Value = 904.8
$ws.Range("K100").Value = 1809.6
$ws.Range("M100").Value = -1268.6
$ws.Range("H105").Value = 45250
$ws.Range("J105").Value = 45250
$ws.Range("L105").Value = 45250
$ws.Range("N105").Value = -52238
$ws.Range("H122").Value = 2407.6667
$ws.Range("I122").Value = 2305
$ws.Range("K122").Value = 6915
$ws.Range("M122").Value = -4465
$ws.Range("H132").Value = 3611.2307
$ws.Range("I132").Value = 3540.5454
$ws.Range("K132").Value = 10621.6362
$ws.Range("M132").Value = -8091.636200000001
$ws.Range("H136").Value = 13267.467
$ws.Range("I136").Value = 5030.636
$ws.Range("K136").Value = 15091.908
$ws.Range("M136").Value = -12541.908
